# Fruta / hortaliza, semanal
# Updates the per-row market data (Fecha, Calidad, Volumen, Precio minimo/maximo/
# promedio ponderado, Precio $/Kg, and occasionally Variedad) for rows 2-39 on the
# active sheet to reflect the refreshed weekly values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44305
$ws.Cells.Item(2, 10).Value = 120
$ws.Cells.Item(2, 11).Value = 13000
$ws.Cells.Item(2, 12).Value = 14000
$ws.Cells.Item(2, 13).Value = 13500
$ws.Cells.Item(2, 16).Value = 900
$ws.Cells.Item(3, 4).Value = 44270
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 120
$ws.Cells.Item(3, 12).Value = 13000
$ws.Cells.Item(3, 13).Value = 12500
$ws.Cells.Item(3, 16).Value = 833
$ws.Cells.Item(4, 4).Value = 44389
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 21000
$ws.Cells.Item(4, 12).Value = 22000
$ws.Cells.Item(4, 13).Value = 21500
$ws.Cells.Item(4, 16).Value = 1433
$ws.Cells.Item(5, 4).Value = 44225
$ws.Cells.Item(5, 10).Value = 120
$ws.Cells.Item(5, 11).Value = 11000
$ws.Cells.Item(5, 12).Value = 12000
$ws.Cells.Item(5, 13).Value = 11500
$ws.Cells.Item(5, 16).Value = 767
$ws.Cells.Item(6, 4).Value = 44298
$ws.Cells.Item(6, 10).Value = 120
$ws.Cells.Item(6, 11).Value = 16000
$ws.Cells.Item(6, 12).Value = 17000
$ws.Cells.Item(6, 13).Value = 16500
$ws.Cells.Item(6, 16).Value = 1100
$ws.Cells.Item(7, 4).Value = 44232
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 17500
$ws.Cells.Item(7, 16).Value = 1167
$ws.Cells.Item(8, 4).Value = 44232
$ws.Cells.Item(8, 9).Value = "Segunda"
$ws.Cells.Item(8, 11).Value = 15000
$ws.Cells.Item(8, 12).Value = 16000
$ws.Cells.Item(8, 13).Value = 15500
$ws.Cells.Item(8, 16).Value = 1033
$ws.Cells.Item(9, 4).Value = 44333
$ws.Cells.Item(9, 9).Value = "Segunda"
$ws.Cells.Item(9, 11).Value = 24000
$ws.Cells.Item(9, 12).Value = 25000
$ws.Cells.Item(9, 13).Value = 24500
$ws.Cells.Item(9, 16).Value = 1633
$ws.Cells.Item(10, 4).Value = 44162
$ws.Cells.Item(10, 10).Value = 170
$ws.Cells.Item(11, 4).Value = 44162
$ws.Cells.Item(11, 9).Value = "Segunda"
$ws.Cells.Item(11, 10).Value = 200
$ws.Cells.Item(12, 4).Value = 44165
$ws.Cells.Item(12, 10).Value = 160
$ws.Cells.Item(12, 11).Value = 10000
$ws.Cells.Item(12, 12).Value = 11000
$ws.Cells.Item(12, 13).Value = 10500
$ws.Cells.Item(12, 16).Value = 700
$ws.Cells.Item(13, 4).Value = 44165
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 120
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 9000
$ws.Cells.Item(13, 13).Value = 8500
$ws.Cells.Item(13, 16).Value = 567
$ws.Cells.Item(14, 4).Value = 44169
$ws.Cells.Item(14, 10).Value = 160
$ws.Cells.Item(14, 11).Value = 12000
$ws.Cells.Item(14, 12).Value = 13000
$ws.Cells.Item(14, 13).Value = 12500
$ws.Cells.Item(14, 16).Value = 833
$ws.Cells.Item(15, 4).Value = 44169
$ws.Cells.Item(15, 9).Value = "Segunda"
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 10000
$ws.Cells.Item(15, 12).Value = 11000
$ws.Cells.Item(15, 13).Value = 10500
$ws.Cells.Item(15, 16).Value = 700
$ws.Cells.Item(16, 4).Value = 44431
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 120
$ws.Cells.Item(16, 11).Value = 28000
$ws.Cells.Item(16, 12).Value = 30000
$ws.Cells.Item(16, 13).Value = 29000
$ws.Cells.Item(16, 16).Value = 1933
$ws.Cells.Item(17, 4).Value = 44372
$ws.Cells.Item(17, 10).Value = 140
$ws.Cells.Item(17, 11).Value = 16000
$ws.Cells.Item(17, 12).Value = 17000
$ws.Cells.Item(17, 13).Value = 16714
$ws.Cells.Item(17, 16).Value = 1114
$ws.Cells.Item(18, 4).Value = 44372
$ws.Cells.Item(18, 10).Value = 60
$ws.Cells.Item(18, 11).Value = 12000
$ws.Cells.Item(18, 12).Value = 12000
$ws.Cells.Item(18, 13).Value = 12000
$ws.Cells.Item(18, 16).Value = 800
$ws.Cells.Item(19, 4).Value = 44372
$ws.Cells.Item(19, 9).Value = "Tercera"
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 12).Value = 11000
$ws.Cells.Item(19, 13).Value = 11000
$ws.Cells.Item(19, 16).Value = 733
$ws.Cells.Item(20, 4).Value = 44239
$ws.Cells.Item(20, 10).Value = 140
$ws.Cells.Item(20, 11).Value = 16000
$ws.Cells.Item(20, 12).Value = 17000
$ws.Cells.Item(20, 13).Value = 16500
$ws.Cells.Item(20, 16).Value = 1100
$ws.Cells.Item(21, 4).Value = 44312
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 11).Value = 14000
$ws.Cells.Item(21, 12).Value = 15000
$ws.Cells.Item(21, 13).Value = 14500
$ws.Cells.Item(21, 16).Value = 967
$ws.Cells.Item(22, 4).Value = 44284
$ws.Cells.Item(22, 11).Value = 11000
$ws.Cells.Item(22, 12).Value = 12000
$ws.Cells.Item(22, 13).Value = 11500
$ws.Cells.Item(22, 16).Value = 767
$ws.Cells.Item(23, 4).Value = 44330
$ws.Cells.Item(23, 10).Value = 120
$ws.Cells.Item(23, 11).Value = 28000
$ws.Cells.Item(23, 12).Value = 30000
$ws.Cells.Item(23, 13).Value = 29000
$ws.Cells.Item(23, 16).Value = 1933
$ws.Cells.Item(24, 4).Value = 44344
$ws.Cells.Item(24, 8).Value = "Cristal"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 11).Value = 24000
$ws.Cells.Item(24, 12).Value = 25000
$ws.Cells.Item(24, 13).Value = 24500
$ws.Cells.Item(24, 16).Value = 1633
$ws.Cells.Item(25, 4).Value = 44340
$ws.Cells.Item(25, 11).Value = 27000
$ws.Cells.Item(25, 12).Value = 28000
$ws.Cells.Item(25, 13).Value = 27500
$ws.Cells.Item(25, 16).Value = 1833
$ws.Cells.Item(26, 4).Value = 44267
$ws.Cells.Item(26, 10).Value = 160
$ws.Cells.Item(26, 11).Value = 12000
$ws.Cells.Item(26, 12).Value = 13000
$ws.Cells.Item(26, 13).Value = 12500
$ws.Cells.Item(26, 16).Value = 833
$ws.Cells.Item(27, 4).Value = 44295
$ws.Cells.Item(27, 11).Value = 19000
$ws.Cells.Item(27, 12).Value = 20000
$ws.Cells.Item(27, 13).Value = 19500
$ws.Cells.Item(27, 16).Value = 1300
$ws.Cells.Item(28, 4).Value = 44295
$ws.Cells.Item(28, 9).Value = "Segunda"
$ws.Cells.Item(28, 10).Value = 120
$ws.Cells.Item(28, 11).Value = 16000
$ws.Cells.Item(28, 12).Value = 17000
$ws.Cells.Item(28, 13).Value = 16500
$ws.Cells.Item(28, 16).Value = 1100
$ws.Cells.Item(29, 4).Value = 44249
$ws.Cells.Item(29, 10).Value = 120
$ws.Cells.Item(29, 11).Value = 12000
$ws.Cells.Item(29, 12).Value = 13000
$ws.Cells.Item(29, 13).Value = 12500
$ws.Cells.Item(29, 16).Value = 833
$ws.Cells.Item(30, 4).Value = 44309
$ws.Cells.Item(30, 11).Value = 11000
$ws.Cells.Item(30, 12).Value = 12000
$ws.Cells.Item(30, 13).Value = 11500
$ws.Cells.Item(30, 16).Value = 767
$ws.Cells.Item(31, 4).Value = 44358
$ws.Cells.Item(31, 10).Value = 160
$ws.Cells.Item(31, 11).Value = 21000
$ws.Cells.Item(31, 12).Value = 22000
$ws.Cells.Item(31, 13).Value = 21500
$ws.Cells.Item(31, 16).Value = 1433
$ws.Cells.Item(32, 4).Value = 44218
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 11).Value = 22000
$ws.Cells.Item(32, 12).Value = 23000
$ws.Cells.Item(32, 13).Value = 22500
$ws.Cells.Item(32, 16).Value = 1500
$ws.Cells.Item(33, 4).Value = 44316
$ws.Cells.Item(33, 8).Value = "Cristal"
$ws.Cells.Item(33, 9).Value = "Segunda"
$ws.Cells.Item(33, 11).Value = 9000
$ws.Cells.Item(33, 12).Value = 10000
$ws.Cells.Item(33, 13).Value = 9500
$ws.Cells.Item(33, 16).Value = 633
$ws.Cells.Item(34, 4).Value = 44176
$ws.Cells.Item(34, 8).Value = "Inferno"
$ws.Cells.Item(34, 11).Value = 11000
$ws.Cells.Item(34, 12).Value = 12000
$ws.Cells.Item(34, 13).Value = 11500
$ws.Cells.Item(34, 16).Value = 767
$ws.Cells.Item(35, 4).Value = 44176
$ws.Cells.Item(35, 8).Value = "Inferno"
$ws.Cells.Item(35, 10).Value = 160
$ws.Cells.Item(36, 4).Value = 44379
$ws.Cells.Item(36, 10).Value = 100
$ws.Cells.Item(36, 11).Value = 19000
$ws.Cells.Item(36, 12).Value = 20000
$ws.Cells.Item(36, 13).Value = 19400
$ws.Cells.Item(36, 16).Value = 1293
$ws.Cells.Item(37, 4).Value = 44435
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 11).Value = 28000
$ws.Cells.Item(37, 12).Value = 30000
$ws.Cells.Item(37, 13).Value = 29000
$ws.Cells.Item(37, 16).Value = 1933
$ws.Cells.Item(38, 4).Value = 44167
$ws.Cells.Item(38, 10).Value = 140
$ws.Cells.Item(38, 11).Value = 11000
$ws.Cells.Item(38, 12).Value = 12000
$ws.Cells.Item(38, 13).Value = 11500
$ws.Cells.Item(38, 16).Value = 767
$ws.Cells.Item(39, 4).Value = 44399
$ws.Cells.Item(39, 11).Value = 12000
$ws.Cells.Item(39, 12).Value = 13000
$ws.Cells.Item(39, 13).Value = 20500
$ws.Cells.Item(39, 16).Value = 1367
